# Trade #62 closed at 2026-02-17 08:48:31 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up numbers to reflect the newly
# closed trade, and appends the new trade row to both the "All Trades" and
# "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Summary sheet
# ----------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.65
$summary.Range("B4").Value = -0.35
$summary.Range("B5").Value = -0.11
$summary.Range("B6").Value = 62
$summary.Range("B8").Value = 26
$summary.Range("B9").Value = 40.32

# ----------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ----------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.65000000000001
$status.Range("D4").Value = 62
$status.Range("E4").Value = -0.35
$status.Range("F4").Value = -0.35
$status.Range("G4").Value = 40.32

# ----------------------------------------------------------------------
# New trade row appended to "All Trades" and "MarketMaking" sheets
#
# Column B holds a date-shaped string ("2026-02-17"); a leading apostrophe
# forces Excel to keep it as literal text instead of auto-converting it to
# a date serial number (matches the source workbook's inlineStr cells).
# ----------------------------------------------------------------------
$newRow = @(62, "'2026-02-17", "08:48:25", "MarketMaking", "UP", 0.85, 0.8100000000000001, "CLOSED", -4.7059, -0.04, 99.65000000000001, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.14)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 63
    for ($col = 1; $col -le $newRow.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $newRow[$col - 1]
    }
}
